# Adds a new "2022-Q1" sheet (before "总计") with per-fund holding data, and
# updates the "总计" (totals) summary sheet with a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet, positioned right after "2021-Q4"
#    and right before "总计". We duplicate the "2021-Q4" sheet (instead of
#    Worksheets.Add()) so that sheet-level properties (outline/page setup,
#    default row height, etc.) and existing cell styles are inherited
#    rather than re-created from scratch.
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Copy($null, $q4Sheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The duplicated sheet only has 6 data rows (rows 2-7); extend the styled
# index column down to row 15 to cover all 14 data rows.
$newSheet.Range("A2").Copy()
$newSheet.Range("A8:A15").PasteSpecial($xlPasteFormats)

# Header row.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Per-fund holding data for 2022-Q1.
$fundRows = @(
    @{Row=2;  A=0;  B="002624"; C="广发优企精选灵活配置混合A";     D="12.98"; E="92.40"; F="7.97"; G="1.0345"; H=5},
    @{Row=3;  A=1;  B="011866"; C="广发价值增长混合型证券投资基金A"; D="18.06"; E="92.01"; F="5.28"; G="0.9536"; H=8},
    @{Row=4;  A=2;  B="270025"; C="广发行业领先混合A";             D="11.11"; E="91.67"; F="6.70"; G="0.7444"; H=7},
    @{Row=5;  A=3;  B="960001"; C="广发行业领先混合H";             D="11.11"; E="91.67"; F="6.70"; G="0.7444"; H=7},
    @{Row=6;  A=4;  B="160726"; C="嘉实瑞享定期开放灵活配置混合";   D="23.58"; E="63.95"; F="2.40"; G="0.5659"; H=7},
    @{Row=7;  A=5;  B="000747"; C="广发逆向策略灵活配置混合";       D="1.25";  E="89.61"; F="6.83"; G="0.0854"; H=7},
    @{Row=8;  A=6;  B="210002"; C="金鹰红利价值混合";               D="0.88";  E="69.18"; F="8.17"; G="0.0719"; H=1},
    @{Row=9;  A=7;  B="011867"; C="广发价值增长混合型证券投资基金C"; D="0.81";  E="92.01"; F="5.28"; G="0.0428"; H=8},
    @{Row=10; A=8;  B="011765"; C="兴银高端制造混合A";             D="1.01";  E="93.23"; F="3.11"; G="0.0314"; H=3},
    @{Row=11; A=9;  B="000994"; C="建信睿盈灵活配置混合A";         D="0.58";  E="88.02"; F="3.61"; G="0.0209"; H=8},
    @{Row=12; A=10; B="011766"; C="兴银高端制造混合C";             D="0.39";  E="93.23"; F="3.11"; G="0.0121"; H=3},
    @{Row=13; A=11; B="010021"; C="广发优企精选灵活配置混合C";     D="0.15";  E="92.40"; F="7.97"; G="0.0120"; H=5},
    @{Row=14; A=12; B="000995"; C="建信睿盈灵活配置混合C";         D="0.19";  E="88.02"; F="3.61"; G="0.0069"; H=8},
    @{Row=15; A=13; B="005146"; C="兴银丰润灵活配置混合";           D="0.05";  E="93.36"; F="3.45"; G="0.0017"; H=5}
)

foreach ($fr in $fundRows) {
    $r = $fr.Row
    $newSheet.Cells.Item($r, 1).Value = $fr.A
    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $fr.B
    $newSheet.Cells.Item($r, 3).Value = $fr.C
    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $fr.D
    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $fr.E
    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $fr.F
    $newSheet.Cells.Item($r, 7).NumberFormat = "@"
    $newSheet.Cells.Item($r, 7).Value = $fr.G
    $newSheet.Cells.Item($r, 8).Value = $fr.H
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a new leading data row for
#    2022-Q1 and push the existing rows down by one.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Push existing data rows (2..6) down to (3..7), preserving their B/C/D
# values as-is (column A is recomputed below as a fresh 0-based index).
for ($r = 6; $r -ge 2; $r--) {
    $dst = $r + 1
    $bVal = $totalSheet.Cells.Item($r, 2).Value()
    $cVal = $totalSheet.Cells.Item($r, 3).Value()
    $dVal = $totalSheet.Cells.Item($r, 4).Value()
    $totalSheet.Cells.Item($dst, 2).Value = $bVal
    $totalSheet.Cells.Item($dst, 3).Value = $cVal
    $totalSheet.Cells.Item($dst, 4).Value = $dVal
}

# New first data row: 2022-Q1.
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 14
$totalSheet.Cells.Item(2, 4).Value = 4.33

# Make sure column A keeps the same formatting (style) for the
# newly-created row 7, then recompute the sequential index for all rows.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A7").PasteSpecial($xlPasteFormats)

for ($r = 2; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# Restore the originally-active sheet/tab (the first sheet), since copying
# a worksheet makes the copy the active tab.
$wb.Worksheets.Item(1).Activate()

Write-Output "2022-Q1 sheet added; 总计 sheet updated."
